# Add "0x" prefix to every colon-separated hex byte token in the
# "doip" (column G) and "uds" (column H) columns, for every data row.
# Values equal to "N/A" (or otherwise containing no ":" separated hex
# bytes) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    foreach ($col in 7, 8) {
        $cell = $ws.Cells.Item($r, $col)
        $orig = $cell.Value2

        if ($orig -ne $null -and $orig -ne "N/A" -and $orig -ne "") {
            $parts = $orig -split ":"
            $newParts = @()
            foreach ($p in $parts) {
                $newParts += "0x" + $p
            }
            $cell.Value = ($newParts -join ":")
        }
    }
}
